$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:G values per row (r2..r57), derived from regenerated s_vals after filtering save games
$data = @(
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.594575437922795),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642),
    @(0.02258322285507441, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 1, 1.012960547955778),
    @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 5.964442013611383),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(0.02258322285507441, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 1.572787616952587),
    @(0.1554434735375247, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 0, 1.420375733316612),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245),
    @(0.3464964993005633, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1, 1.051601690082842),
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 5.964442013611383),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245),
    @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 6.741336633845642),
    @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 5.964442013611383),
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245),
    @(0.1554434735375247, 0.004309184025731883, 0.7127328510149897, 0.4998867070740569, 0, 1.372372215652303),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729),
    @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.594575437922795),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.594575437922795),
    @(0.001754667048134761, 0.05231270169004087, 16.98373111632243, 0.4998867070740569, 1, 17.53768519213466),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 1.896700893398075),
    @(0.7287194209349384, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 19.86557370323023),
    @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 1.705647867635037),
    @(0.001754667048134761, 0.004309184025731883, 3.082599426703578, 0.4998867070740569, 0, 3.588549984851502),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.3464964993005633, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 2.652525447291612),
    @(0.006876353814593728, 0.05231270169004087, 3.082599426703578, 0.4998867070740569, 0, 3.64167518928227),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126),
)

$arr = New-Object 'object[,]' 56,6
for ($i = 0; $i -lt $data.Count; $i++) {
    for ($j = 0; $j -lt 6; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$ws.Range("B2:G57").Value = $arr
